# "Lista de riscos.docx" - teste 2 - Documento de Acompanhamento
#
# Changes applied (per commit diff):
#   1. Remove the paragraph "Falta de recursos para cobrir os custos".
#   2. Replace the text "Troca do computador do cliente enquanto o
#      software é desenvolvido" with "Erros não detectados na fase de
#      testes" (the wording that used to live in the final paragraph).
#   3. Remove the now-duplicated trailing paragraph that used to hold
#      "Erros não detectados na fase de testes" (together with the
#      stray _GoBack bookmark it carried).
#   4. Header/footer distance collapses to 0 (was 708 twips ~ 35.4pt).

$d = $word.ActiveDocument

function Get-ParagraphIndex($doc, $prefix) {
    $i = 0
    $found = -1
    foreach ($p in $doc.Paragraphs) {
        $i = $i + 1
        if ($p.Range.Text -like "$prefix*") {
            $found = $i
        }
    }
    return $found
}

# 1. Delete the "Falta de recursos para cobrir os custos" paragraph entirely.
$iRecursos = Get-ParagraphIndex $d "Falta de recursos para cobrir os custos"
if ($iRecursos -gt 0) {
    $d.Paragraphs($iRecursos).Range.Delete()
}

# 2. Update the "Troca do computador..." paragraph to the new wording.
$iTroca = Get-ParagraphIndex $d "Troca do computador do cliente"
if ($iTroca -gt 0) {
    $d.Paragraphs($iTroca).Range.Text = "Erros não detectados na fase de testes"
}

# 3. Drop the old trailing paragraph (duplicate text + _GoBack bookmark).
$lastIndex = $d.Paragraphs.Count
$d.Paragraphs($lastIndex).Range.Delete()

# 4. Collapse header/footer distance to 0.
$section = $d.Sections(1)
$section.PageSetup.HeaderDistance = 0
$section.PageSetup.FooterDistance = 0
